$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.390.07"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "'2.069.93"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'235.22"
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'57.46"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "'0.393"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "'2.373.04"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'14.45"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "'20.83"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'0.777"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'5.18"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "'2.067.46"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'37.334.95"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'69.57"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'0.0₃0819"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "'226.96"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "'167.23"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'1.41"
$ws.Range("E28").Value = "  -6.79%  "
$ws.Range("D29").Value = "'0.129"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").Value = "'19.11"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "'0.117"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "'4.54"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'2.46"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'3.35"
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.78"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("D42").Value = "'1.486.55"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "'97.69"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'0.0213"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "'4.19"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'15.12"
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.20"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "'47.64"
$ws.Range("E51").Value = "  +6.20%  "
